$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K13").Value = "['Draw', 0, 0]"
$ws.Range("L13").Value = "['Belgium', 2, 1]"
$ws.Range("M13").Value = "['France', 0, 3]"
$ws.Range("N13").Value = "['Turkey', 2, 1]"
$ws.Range("O13").Value = "['Portugal', 2, 1]"
$ws.Range("P13").Value = "['Croatia', 2, 0]"
$ws.Range("Q13").Value = "['Germany', 3, 1]"
